# This script reproduces the "weekly" update described in the commit message
# "Fruta / hortaliza, semanal": a new week of price data (dated 2022-05-20,
# serial 44701) is inserted as two new rows (Primera / Segunda quality) right
# after the header's first existing data block, pushing every subsequent row
# down by two positions (old row N -> new row N+2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 16 (this shifts rows 16:97 down to 18:99,
# and Excel copies the formatting - including the date number format on
# column D - from the row above, matching the original file's style).
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# New row 16: "Primera" quality entry for the new week
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 44701
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100112036
$ws.Cells.Item(16, 7).Value = "Caigua"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 120
$ws.Cells.Item(16, 11).Value = 14000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 14500
$ws.Cells.Item(16, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(16, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 16).Value = 725
$ws.Cells.Item(16, 17).Value = 20
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# New row 17: "Segunda" quality entry for the new week
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(17, 4).Value = 44701
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = 100112036
$ws.Cells.Item(17, 7).Value = "Caigua"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Segunda"
$ws.Cells.Item(17, 10).Value = 140
$ws.Cells.Item(17, 11).Value = 11000
$ws.Cells.Item(17, 12).Value = 12000
$ws.Cells.Item(17, 13).Value = 11500
$ws.Cells.Item(17, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(17, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value = 575
$ws.Cells.Item(17, 17).Value = 20
$ws.Cells.Item(17, 18).Value = "Hortaliza"
